# Daily attendance processing - swap "Recorded By" name ordering in column G
# so that entries formatted as "<email>, System" become "System, <email>".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count

for ($r = 1; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)  # Column G = 7
    $val = $cell.Text

    if ($val -ne $null -and $val -match '^(dnasr281@gmail\.com|admin@admin\.com), System$') {
        $email = $matches[1]
        $cell.Value = "System, $email"
    }
}
